$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 185 (existing rows 185-258 shift down to 186-259),
# adding a new "z" / "Depth/height (meters)" / "Profondeur/hauteur (metres)"
# translation entry (aggregation_type normalization).
$ws.Rows.Item(185).Insert()

$ws.Range("A185").Value = "z"
$ws.Range("C185").Value = "Depth/height (meters)"
$ws.Range("D185").Value = "Profondeur/hauteur (mètres)"

# Match the author's final cursor position/selection.
$ws.Range("D185").Select()
